$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Mobile Number")

$numbers = @(
    "9481758105",
    "2718454523",
    "2154478176",
    "3728940729",
    "7261910232",
    "8105886786",
    "8647284553",
    "5373968912",
    "4131060765",
    "4087557155",
    "6582291831",
    "4956187639",
    "1604155166",
    "4445460072",
    "0120437413",
    "1984094441",
    "8372006150",
    "9737934498",
    "4462166045",
    "7472031041",
    "2272891346",
    "4764757375",
    "9618079294",
    "6746417643",
    "6358885530",
    "3942262542",
    "5829839274",
    "2464542431",
    "4936201314",
    "7521146830",
    "1031496065",
    "3822327461",
    "9559405562",
    "3040277876",
    "7514424092",
    "6770262251",
    "0939873681",
    "8356831008",
    "9294548500",
    "2752916912",
    "3573140489",
    "5768779175",
    "6263099966",
    "1852755047",
    "2787555746",
    "8637837862",
    "6517517021",
    "5826411890",
    "9712625243",
    "2023262132",
    "9803240222",
    "9965146585",
    "9141897296",
    "4539727291",
    "2264992038",
    "4688608568",
    "2042397159",
    "1740477750",
    "8331408592",
    "0519256949",
    "6987426695",
    "8455189902",
    "1934961017",
    "5104174812",
    "1919940069",
    "3919811534",
    "2924347501",
    "1452040359",
    "4902488420",
    "9388864286",
    "9057072236",
    "9798538334",
    "1351428444",
    "1036775979",
    "4571959374",
    "1059585427",
    "9100554479",
    "3674500568",
    "6118490504",
    "7391918471",
    "2221750282",
    "3033728593",
    "0107710010",
    "3302772706",
    "7581740116",
    "4711103563",
    "6351813634",
    "3926414576",
    "4845749618"
)

$startRow = 29
$endRow = $startRow + $numbers.Length - 1
$rng = $ws.Range("A" + $startRow + ":A" + $endRow)

# Force text interpretation so values with leading zeros keep them,
# matching the existing phone-number rows above.
$rng.NumberFormat = "@"

for ($i = 0; $i -lt $numbers.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A" + $row).Value = $numbers[$i]
}

# Restore the default "Normal" style so the new cells match the
# unstyled (s=0) cells already on the sheet.
$rng.Style = "Normal"
